$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# D-column (Price) cells need NumberFormat forced to text so Excel doesn't
# reinterpret numeric-looking strings (e.g. '1.00', '516.74') as numbers and
# drop formatting. ClearFormats() afterwards removes the now-unneeded style
# index so the cell's style stays identical to the original (unstyled).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.354.48'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +1.51%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.358.09'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +1.50%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '516.74'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '135.57'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +1.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.24%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.539'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.77%  '
$ws.Range("E9").Value = '  -0.18%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.45'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +5.04%  '
$ws.Range("E11").Value = '  -0.85%  '
$ws.Range("E12").Value = '  +0.77%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '24.43'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +2.88%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.779.53'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.38%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '57.366.42'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.26%  '
$ws.Range("E16").Value = '  +0.40%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.359.59'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +1.46%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.53'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.74%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '328.27'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +2.12%  '
$ws.Range("E20").Value = '  -0.23%  '
$ws.Range("E21").Value = '  +0.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '61.19'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.73%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.80'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +13.24%  '
$ws.Range("E25").Value = '  +4.45%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.995'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.23%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.32'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +9.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0741'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.51%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '167.05'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -2.44%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.69'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.70%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.25'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.33%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.55'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +1.55%  '
$ws.Range("E33").Value = '  +0.03%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.29'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +2.74%  '
$ws.Range("E35").Value = '  +0.27%  '
$ws.Range("E36").Value = '  +0.95%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.914'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -4.71%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.59'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +4.83%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '38.98'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +4.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '149.64'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +7.42%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.387'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +1.57%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.65'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +2.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '287.27'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +3.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.35'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +4.72%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0935'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.60%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0507'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.564'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.79%  '
$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '18.19'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +5.70%  '
$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0218'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +1.29%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.53'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +3.85%  '
$ws.Range("B51").Value = 'WhiteBITCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '11.00'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +1.41%  '
